$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ministry Course Code and Level" (column G, e.g. "ENST 12") is being split
# into two separate columns: "Ministry Course Code" (G, e.g. "ENST") and a
# newly inserted "Ministry Course Level" (H, e.g. 12). Inserting a column at
# H shifts Session Date / Final Percent / Final Letter Grade / Credits one
# column to the right (H->I, I->J, J->K, K->L).
$ws.Columns.Item(8).Insert()

# Header row
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("G2").Value = "ENST"
$ws.Range("G3").Value = "ENST"
$ws.Range("G4").Value = "ENST"
$ws.Range("H1").Value = "Ministry Course Level"
$ws.Range("H2").Value = 12
$ws.Range("H3").Value = 12
$ws.Range("H4").Value = 12

$ws.Range("G1:H1048576").Select() | Out-Null
